$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:H1) so the plain numeric labels become the
# "scatterN" text labels used to identify each data series.
$ws.Range("A1").Value = "scatter1"
$ws.Range("B1").Value = "scatter2"
$ws.Range("C1").Value = "scatter3"
$ws.Range("D1").Value = "scatter4"
$ws.Range("E1").Value = "scatter5"
$ws.Range("F1").Value = "scatter6"
$ws.Range("G1").Value = "scatter7"
$ws.Range("H1").Value = "scatter8"

# Reflect the refreshed selection/active cell position.
$ws.Range("J5").Select()
